$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9.684777666666667
$ws.Range("H2").Value = 29.054333
$ws.Range("I2").Value = 0.8593785179720864
$ws.Range("J2").Value = 0.8593785179720863
$ws.Range("M2").Value = 25.37147633333333
$ws.Range("N2").Value = 76.114429
$ws.Range("O2").Value = 0.5780881462719274
$ws.Range("P2").Value = 0.5780881462719274
$ws.Range("Q2").Value = 245.7171073634286
$ws.Range("R2").Value = 2211.453966270857
$ws.Range("S2").Value = 0.4967965344003997
$ws.Range("T2").Value = 0.4967965344003996

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 9.684777666666667
$ws.Range("H3").Value = 29.054333
$ws.Range("I3").Value = 0.8593785179720864
$ws.Range("J3").Value = 0.8593785179720863
$ws.Range("O3").Value = 0.2328552951919536
$ws.Range("P3").Value = 0.2328552951919536
$ws.Range("Q3").Value = 98.97544161355279
$ws.Range("R3").Value = 890.778974521975
$ws.Range("S3").Value = 0.2001108384840138
$ws.Range("T3").Value = 0.2001108384840138

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 9.684777666666667
$ws.Range("H4").Value = 29.054333
$ws.Range("I4").Value = 0.8593785179720864
$ws.Range("J4").Value = 0.8593785179720863
$ws.Range("O4").Value = 0.189056558536119
$ws.Range("P4").Value = 0.189056558536119
$ws.Range("Q4").Value = 80.35873247213777
$ws.Range("R4").Value = 723.22859224924
$ws.Range("S4").Value = 0.1624711450876729
$ws.Range("T4").Value = 0.1624711450876729

# Row 5
$ws.Range("I5").Value = 0.05551925637723486
$ws.Range("J5").Value = 0.05551925637723486
$ws.Range("M5").Value = 25.37147633333333
$ws.Range("N5").Value = 76.114429
$ws.Range("O5").Value = 0.5780881462719274
$ws.Range("P5").Value = 0.5780881462719274
$ws.Range("Q5").Value = 15.87429845485833
$ws.Range("R5").Value = 142.868686093725
$ws.Range("S5").Value = 0.03209502400151159
$ws.Range("T5").Value = 0.03209502400151159

# Row 6
$ws.Range("I6").Value = 0.05551925637723486
$ws.Range("J6").Value = 0.05551925637723486
$ws.Range("O6").Value = 0.2328552951919536
$ws.Range("P6").Value = 0.2328552951919536
$ws.Range("S6").Value = 0.01292795283255878
$ws.Range("T6").Value = 0.01292795283255878

# Row 7
$ws.Range("I7").Value = 0.05551925637723486
$ws.Range("J7").Value = 0.05551925637723486
$ws.Range("O7").Value = 0.189056558536119
$ws.Range("P7").Value = 0.189056558536119
$ws.Range("S7").Value = 0.0104962795431645
$ws.Range("T7").Value = 0.0104962795431645

# Row 8
$ws.Range("I8").Value = 0.0851022256506788
$ws.Range("J8").Value = 0.08510222565067879
$ws.Range("M8").Value = 25.37147633333333
$ws.Range("N8").Value = 76.114429
$ws.Range("O8").Value = 0.5780881462719274
$ws.Range("P8").Value = 0.5780881462719274
$ws.Range("Q8").Value = 24.33278500656422
$ws.Range("R8").Value = 218.995065059078
$ws.Range("S8").Value = 0.04919658787001618
$ws.Range("T8").Value = 0.04919658787001617

# Row 9
$ws.Range("I9").Value = 0.0851022256506788
$ws.Range("J9").Value = 0.08510222565067879
$ws.Range("O9").Value = 0.2328552951919536
$ws.Range("P9").Value = 0.2328552951919536
$ws.Range("S9").Value = 0.01981650387538106
$ws.Range("T9").Value = 0.01981650387538106

# Row 10
$ws.Range("I10").Value = 0.0851022256506788
$ws.Range("J10").Value = 0.08510222565067879
$ws.Range("O10").Value = 0.189056558536119
$ws.Range("P10").Value = 0.189056558536119
$ws.Range("S10").Value = 0.01608913390528156
$ws.Range("T10").Value = 0.01608913390528156
